$d = $word.ActiveDocument

# --- Locate the two target paragraphs by their (original) text ---
$targetOld = $null
$targetNew = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Will also store frequency of clues remaining in a hash table*") {
        $targetOld = $p
    }
    if ($p.Range.Text -like "When creating the structure to hold the remaining cells*") {
        $targetNew = $p
    }
}

if ($targetOld -eq $null) {
    throw "Could not find paragraph starting with 'Will also store frequency...'"
}
if ($targetNew -eq $null) {
    throw "Could not find paragraph starting with 'When creating the structure...'"
}

# --- Paragraph 1: rewrite "Will also store..." paragraph with expanded text ---
$xmlPara1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>Will also store frequency of clues remaining in a</w:t></w:r>
<w:r><w:t>n array</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
<w:r><w:t xml:space="preserve"> Index is key-1 (</w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>ie</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> for grid size 9, keys are numbers 1-9 minus 1 to obtain index). Value</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:r><w:t xml:space="preserve">stored will be the frequency of numbers remaining. </w:t></w:r>
<w:r><w:t xml:space="preserve">This can be used to ensure that at least </w:t></w:r>
<w:r><w:t xml:space="preserve">1 clue remains for </w:t></w:r>
<w:r><w:t>SIZE-1</w:t></w:r>
<w:r><w:t xml:space="preserve"> (</w:t></w:r>
<w:r><w:t>all except 1 number</w:t></w:r>
<w:r><w:t>)</w:t></w:r>
<w:r><w:t xml:space="preserve"> at all times</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$targetOld.Range.InsertXML($xmlPara1)

# --- Paragraph 2: replace "When creating the structure..." paragraph with the
#     new "Rem_grid_nums" notes (4 paragraphs total) ---
$xmlPara2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Rem_grid_nums</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> potential slow performance areas:</w:t></w:r>
</w:p>
<w:p>
<w:r><w:tab/><w:t>Updates access node &#8211; comparison check on each number removed</w:t></w:r>
</w:p>
<w:p>
<w:r><w:tab/></w:r>
</w:p>
<w:p/>
<w:p/>
<w:p/>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$targetNew.Range.InsertXML($xmlPara2)
